# Refresh the crypto symbol/price table ("Updated symbol list ... with GitHub
# Actions"). All cells in this sheet are stored as text, so every write below
# is prefixed with a leading apostrophe to force Excel to keep it as a text
# value instead of auto-coercing numeric- or percentage-looking strings
# (e.g. "291.11", "-5.87%", "14") into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every data row (2-51) has its "Hora" (G) value bumped from 13 to 14
$ws.Range("G2:G51").Value = "'14"

# Per-row updates for Coin/Link/Price/Volume columns
# Row 2
$ws.Range("D2").Value = "'291.11"
$ws.Range("E2").Value = "'-5.87%"
# Row 3
$ws.Range("D3").Value = "'39.79"
$ws.Range("E3").Value = "'-2.97%"
# Row 4
$ws.Range("D4").Value = "'5.031"
$ws.Range("E4").Value = "'-2.98%"
# Row 5
$ws.Range("D5").Value = "'0.07346"
$ws.Range("E5").Value = "'-4.05%"
# Row 6
$ws.Range("B6").Value = "'GateToken"
$ws.Range("C6").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.289"
$ws.Range("E6").Value = "'-0.34%"
# Row 7
$ws.Range("B7").Value = "'FTXToken"
$ws.Range("C7").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.554"
$ws.Range("E7").Value = "'-10.67%"
# Row 8
$ws.Range("B8").Value = "'MXToken"
$ws.Range("C8").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9166"
$ws.Range("E8").Value = "'0.08%"
# Row 9
$ws.Range("B9").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1192"
$ws.Range("E9").Value = "'-4.51%"
# Row 10
$ws.Range("B10").Value = "'WazirX"
$ws.Range("C10").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1724"
$ws.Range("E10").Value = "'-5.61%"
# Row 11
$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08770"
$ws.Range("E11").Value = "'-4.15%"
# Row 12
$ws.Range("B12").Value = "'BitrueCoin"
$ws.Range("C12").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04183"
$ws.Range("E12").Value = "'1.12%"
# Row 13
$ws.Range("B13").Value = "'BitMartToken"
$ws.Range("C13").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1053"
$ws.Range("E13").Value = "'0.24%"
# Row 14
$ws.Range("B14").Value = "'BitForexToken"
$ws.Range("C14").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001279"
$ws.Range("E14").Value = "'-0.28%"
# Row 15
$ws.Range("B15").Value = "'TigerCash"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005838"
$ws.Range("E15").Value = "'1.23%"
# Row 16
$ws.Range("B16").Value = "'LEO"
$ws.Range("C16").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.403"
$ws.Range("E16").Value = "'1.44%"
# Row 18
$ws.Range("D18").Value = "'0.3296"
$ws.Range("E18").Value = "'-0.71%"
# Row 19
$ws.Range("D19").Value = "'7.594"
$ws.Range("E19").Value = "'3.05%"
# Row 20
$ws.Range("D20").Value = "'0.1354"
$ws.Range("E20").Value = "'-0.10%"
# Row 21
$ws.Range("D21").Value = "'0.2880"
$ws.Range("E21").Value = "'5.62%"
# Row 22
$ws.Range("D22").Value = "'0.03843"
$ws.Range("E22").Value = "'-4.40%"
# Row 23
$ws.Range("D23").Value = "'0.001280"
$ws.Range("E23").Value = "'0.86%"
# Row 24
$ws.Range("D24").Value = "'0.003876"
$ws.Range("E24").Value = "'-5.46%"
# Row 25
$ws.Range("D25").Value = "'0.0001282"
$ws.Range("E25").Value = "'-1.60%"
# Row 26
$ws.Range("D26").Value = "'0.0003733"
# Row 38
$ws.Range("D38").Value = "'0.02326"
$ws.Range("E38").Value = "'-7.99%"
# Row 39
$ws.Range("D39").Value = "'0.05014"
$ws.Range("E39").Value = "'-5.32%"
# Row 40
$ws.Range("D40").Value = "'0.007683"
$ws.Range("E40").Value = "'-2.18%"
# Row 41
$ws.Range("E41").Value = "'172.56%"
# Row 42
$ws.Range("D42").Value = "'0.1268"
$ws.Range("E42").Value = "'-3.04%"
# Row 43
$ws.Range("D43").Value = "'0.007380"
$ws.Range("E43").Value = "'10.95%"
# Row 44
$ws.Range("D44").Value = "'0.007704"
$ws.Range("E44").Value = "'-5.41%"
# Row 45
$ws.Range("D45").Value = "'0.3149"
$ws.Range("E45").Value = "'2.72%"
# Row 46
$ws.Range("D46").Value = "'0.00006525"
$ws.Range("E46").Value = "'-4.00%"
# Row 47
$ws.Range("E47").Value = "'0.03%"
# Row 48
$ws.Range("E48").Value = "'7.64%"
# Row 49
$ws.Range("D49").Value = "'0.004212"
$ws.Range("E49").Value = "'35.71%"
# Row 50
$ws.Range("D50").Value = "'0.00002106"
$ws.Range("E50").Value = "'0.03%"
# Row 51
$ws.Range("D51").Value = "'0.0002006"
$ws.Range("E51").Value = "'0.03%"
